$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 89 and 90 need their data (columns B, F:AC) swapped.
$cols = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

$row89vals = @{}
$row90vals = @{}

foreach ($c in $cols) {
    $row89vals[$c] = $ws.Range("$c" + "89").Value2
    $row90vals[$c] = $ws.Range("$c" + "90").Value2
}

foreach ($c in $cols) {
    $ws.Range("$c" + "89").Value2 = $row90vals[$c]
    $ws.Range("$c" + "90").Value2 = $row89vals[$c]
}
